$wb = $excel.ActiveWorkbook

# Locate the two sheets by their known names (order-independent lookup).
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")

# Reorder tabs so review_info comes first, hotel_info second.
$reviewInfo.Move($hotelInfo)

# Worksheet object references above are position-bound, so after the Move()
# re-resolve the hotel_info sheet by name to get a handle to the right tab.
$hotelInfo = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column into hotel_info right before the City column
# (City was column C, so inserting at column C pushes City -> D, etc.)
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"
